$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same event rows and need the
# same update to column F (想去人数 / interested-count) for rows 2 and 3.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 147
    $ws.Range("F3").Value = 101
}
